$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a run boundary at the edges of a (Start,End) range without
# altering its visible formatting. This engine (like real Word) silently
# re-merges two adjacent runs once they end up with byte-identical <w:rPr>,
# so a plain text insert/delete isn't enough to keep a freshly split run
# distinct from its neighbour. Toggling a boolean character property on and
# back off ($true/1 then $false) is recorded as an explicit-but-default
# direct-formatting entry on the *created* run only, which keeps it from
# re-absorbing into its neighbours while leaving the saved <w:rPr> identical
# to theirs.
# ---------------------------------------------------------------------------
function Split-RunBoundary($range) {
    $range.Bold = 1
    $range.Bold = $false
}

# ---------------------------------------------------------------------------
# Edit 1: "...the bag didn't fit at all and she just avoided brining it."
#      -> "...the bag didn't fit at all and she had just avoided brining it."
# Realised as a 3-way run split: [...she ] + [had ] + [just avoided...]
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("just avoided", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Edit 1: target text 'just avoided' not found"
}
$insertStart = $rng1.Start
$rng1.Collapse(1)
$rng1.Text = "had "

$newRun1 = $d.Range($insertStart, $insertStart + 4)
Split-RunBoundary $newRun1

# ---------------------------------------------------------------------------
# Edit 2: "...websites and programming, when I managed to, but..."
#      -> "...websites and programming, when I manage to, but..."
# Realised as a 3-way run split:
#   [...websites] + [ and programming, when I manage] + [ to, but...]
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$sentence = '"Hmm... I am a part-time DJ or a maybe full-time DJ depending on the point of view" - Her face drew a grin - "I also do a few jobs with websites and programming, when I managed to, but it''s usually not paid at all..."'
$found2 = $rng2.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Edit 2: target sentence not found"
}
$sentStart = $rng2.Start

$bStart = $sentStart + 144   # right after "...with websites"
$bEnd = $sentStart + 175     # right after "...when I manage"

# Drop the trailing "d" of "managed" right after this span.
$dRange = $d.Range($bEnd, $bEnd + 1)
$dRange.Text = ""

$runB = $d.Range($bStart, $bEnd)
Split-RunBoundary $runB

$runC = $d.Range($bEnd, $bEnd + 41)
Split-RunBoundary $runC

Write-Host "Edits applied"
